$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the logged hours
$ws.Range("C2").Value = 5
$ws.Range("B7").Value = 2
$ws.Range("B8").Value = 3

# Update the selected cell to C3 (matches sheetView selection in diff)
$ws.Range("C3").Select()

$wb.Save()
